# GSC export "Breadcrumbs" workbook: roll the 89-day window forward by one day.
#
# Effect (see commit "updated main GSC export data"):
#   - Chart!A2:A90 holds one text date per row ("yyyy-MM-dd" strings, not real
#     Excel dates). Every date advances by one day: the oldest day
#     (2025-10-06) drops off the front and a new day (2026-01-03) appears at
#     the end, with every other row showing what used to be the next row's
#     date.
#   - Chart!C4:C90 holds the cumulative "Items" count lined up with those
#     dates. Because the calendar date in each row advanced by one day, each
#     row's count becomes whatever count used to belong to the next row
#     (C4 <- old C5, C5 <- old C6, ... C89 <- old C90); the newest/last row
#     (C90) keeps the figure it already had since no further data exists yet.
#   - Column B ("Invalid") is untouched - it is all zeros before and after.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# --- Column C: shift the cumulative counts up by one row -------------------
# Snapshot the current values first (rows 4..90), then write row r the value
# that used to live in row r+1, for r = 4..89. Row 90 is left exactly as-is.
$cVals = @()
for ($r = 4; $r -le 90; $r++) {
    $cVals += ,$ws.Cells.Item($r, 3).Value2
}
for ($r = 4; $r -le 89; $r++) {
    $ws.Cells.Item($r, 3).Value = $cVals[$r - 3]
}

# --- Column A: advance every date by one day --------------------------------
# These cells are plain text (shared strings), not Excel date serials, so we
# force Text formatting before writing and clear it again afterwards - this
# keeps Excel from "helpfully" reinterpreting "2025-10-07" as a date serial
# and leaves the cell's style untouched.
$dateRange = $ws.Range("A2:A90")
$dateRange.NumberFormat = "@"
$base = Get-Date -Year 2025 -Month 10 -Day 6
for ($r = 2; $r -le 90; $r++) {
    $d = $base.AddDays($r - 2 + 1)
    $ws.Cells.Item($r, 1).Value = $d.ToString("yyyy-MM-dd")
}
$dateRange.ClearFormats()
